# Daily attendance processing - 2025-12-06 17:23:57
# Reverses the order of the comma-separated "Recorded By" names/emails
# in column G for every data row that has more than one entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ","
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        $revParts = $parts[($parts.Length - 1)..0]
        $newVal = [string]::Join(", ", $revParts)
        $cell.Value2 = $newVal
    }
}
